# B6-PowerPoint.pptx edit: 30 Jul 2020
#
# 1) Three tables (on slides 14, 15 and 16) are re-styled from the
#    deck's custom "Table_0" style to the built-in table style
#    {6F710AA6-208A-411C-95C5-B44E8659CC44}.
# 2) The deck's colour theme is switched from the custom "Integral /
#    Red Violet" palette over to the stock "Office" palette (the
#    swap that, at the OOXML level, exchanges the contents of
#    theme1.xml and theme2.xml). The colour values are applied through
#    ThemeColorScheme, which is the only theme surface this host
#    exposes for writes.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------
$newTableStyleId = "{6F710AA6-208A-411C-95C5-B44E8659CC44}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Swap the presentation's colour theme for the stock Office one
$officeColors = @{
    1  = 0x000000   # dk1      000000
    2  = 0xFFFFFF   # lt1      FFFFFF
    3  = 0x6A5444   # dk2      44546A (stored as BGR for RGB())
    4  = 0xE6E6E7   # lt2      E7E6E6
    5  = 0xD59B5B   # accent1  5B9BD5
    6  = 0x317DED   # accent2  ED7D31
    7  = 0xA5A5A5   # accent3  A5A5A5
    8  = 0x00C0FF   # accent4  FFC000
    9  = 0xC47244   # accent5  4472C4
    10 = 0x47AD70   # accent6  70AD47
    11 = 0xC16305   # hlink    0563C1
    12 = 0x724F95   # folHlink 954F72
}

$masterTheme = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($slot in $officeColors.Keys) {
    $masterTheme.Item($slot).RGB = $officeColors[$slot]
}
